# Apply BoM/Costs reference-designator renumbering + timestamp update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('BoM')
$ws.Range('D9').Value = 'C23 C26'
$ws.Range('D10').Value = 'C32 C34'
$ws.Range('D11').Value = 'C31'
$ws.Range('D12').Value = 'C1 C11 C12 C13 C14 C16 C18 C19 C20 C21 C22 C25 C27 C28 C36 C39'
$ws.Range('D13').Value = 'C15 C17'
$ws.Range('D14').Value = 'C2 C3 C4 C5 C6 C7 C8 C9 C10 C24 C29 C30 C33 C35 C38'
$ws.Range('D15').Value = 'C37'
$ws.Range('D17').Value = 'D9'
$ws.Range('D18').Value = 'D8'
$ws.Range('D19').Value = 'D2 D3 D4 D6'
$ws.Range('D21').Value = 'D5 D7'
$ws.Range('D23').Value = 'J7'
$ws.Range('D27').Value = 'J1 J3'
$ws.Range('D28').Value = 'J5 J9 J16 J17 J18 J20'
$ws.Range('Q28').Value = 'EXP1 EXP2 Audio IN Right Audio OUT Left Audio OUT Right Audio IN Left'
$ws.Range('D32').Value = 'L5'
$ws.Range('D34').Value = 'R2 R12'
$ws.Range('D35').Value = 'R8 R9'
$ws.Range('D36').Value = 'R3'
$ws.Range('D37').Value = 'R1 R6'
$ws.Range('D38').Value = 'R4 R5 R10 R11 R16 R17'
$ws.Range('D39').Value = 'R14'
$ws.Range('D40').Value = 'R13'
$ws.Range('D41').Value = 'R15'
$ws.Range('D45').Value = 'U6'

$ws = $wb.Worksheets.Item('DNF')
$ws.Range('D9').Value = 'J2 J4'
$ws.Range('D12').Value = 'J6'
$ws.Range('D13').Value = 'R7'

$ws = $wb.Worksheets.Item('Costs')
$ws.Range('A10').Value = 'C23 C26'
$ws.Range('A11').Value = 'C32 C34'
$ws.Range('A12').Value = 'C31'
$ws.Range('A13').Value = 'C1 C11 C12 C13 C14 C16 C18 C19 C20 C21 C22 C25 C27 C28 C36 C39'
$ws.Range('A14').Value = 'C15 C17'
$ws.Range('A15').Value = 'C2 C3 C4 C5 C6 C7 C8 C9 C10 C24 C29 C30 C33 C35 C38'
$ws.Range('A16').Value = 'C37'
$ws.Range('A18').Value = 'D9'
$ws.Range('A19').Value = 'D8'
$ws.Range('A20').Value = 'D2 D3 D4 D6'
$ws.Range('A22').Value = 'D5 D7'
$ws.Range('A24').Value = 'J7'
$ws.Range('A28').Value = 'J1 J3'
$ws.Range('A29').Value = 'J5 J9 J16 J17 J18 J20'
$ws.Range('A33').Value = 'L5'
$ws.Range('A35').Value = 'R2 R12'
$ws.Range('A36').Value = 'R8 R9'
$ws.Range('A37').Value = 'R3'
$ws.Range('A38').Value = 'R1 R6'
$ws.Range('A39').Value = 'R4 R5 R10 R11 R16 R17'
$ws.Range('A40').Value = 'R14'
$ws.Range('A41').Value = 'R13'
$ws.Range('A42').Value = 'R15'
$ws.Range('A46').Value = 'U6'
$ws.Range('B54').Value = '2023-11-29 20:54:05'

$ws = $wb.Worksheets.Item('Costs (DNF)')
$ws.Range('A10').Value = 'J2 J4'
$ws.Range('A13').Value = 'J6'
$ws.Range('A14').Value = 'R7'
$ws.Range('B17').Value = '2023-11-29 20:54:05'
